$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G8:H8").ClearFormats()
$ws.Range("G8").Value = 0.0
$ws.Range("H8").Value = 0.0

$ws.Range("G9:H9").ClearFormats()
$ws.Range("G9").Value = 0.0
$ws.Range("H9").Value = 0.0

$ws.Range("A10:B10").ClearFormats()
$ws.Range("A10").Value = 1299800.0
$ws.Range("B10").Value = 1780816.0

$ws.Range("G10:H10").ClearFormats()
$ws.Range("G10").Value = 0.0
$ws.Range("H10").Value = 0.0

$ws.Range("G11:H11").ClearFormats()
$ws.Range("G11").Value = 0.0
$ws.Range("H11").Value = 0.0

$ws.Range("G12:H12").ClearFormats()
$ws.Range("G12").Value = 0.0
$ws.Range("H12").Value = 0.0

$ws.Range("G13:H13").ClearFormats()
$ws.Range("G13").Value = 0.0
$ws.Range("H13").Value = 0.0

$ws.Range("G14:H14").ClearFormats()
$ws.Range("G14").Value = 0.0
$ws.Range("H14").Value = 0.0

$ws.Range("G15:H15").ClearFormats()
$ws.Range("G15").Value = -1334448.0
$ws.Range("H15").Value = 600.0

$ws.Range("G16:H16").ClearFormats()
$ws.Range("G16").Value = 0.0
$ws.Range("H16").Value = 0.0

$ws.Range("G17:H17").ClearFormats()
$ws.Range("G17").Value = 0.0
$ws.Range("H17").Value = 0.0

$ws.Range("G18:H18").ClearFormats()
$ws.Range("G18").Value = 0.0
$ws.Range("H18").Value = 0.0

$ws.Range("C19:D19").ClearFormats()
$ws.Range("C19").Value = 0.0
$ws.Range("D19").Value = 1331736.0

$ws.Range("C20:D20").ClearFormats()
$ws.Range("C20").Value = 0.0
$ws.Range("D20").Value = 0.0

$ws.Range("K21:L21").ClearFormats()
$ws.Range("K21").Value = 0.0
$ws.Range("L21").Value = 0.0

$ws.Range("K22:L22").ClearFormats()
$ws.Range("K22").Value = 0.0
$ws.Range("L22").Value = 0.0

$ws.Range("K23:L23").ClearFormats()
$ws.Range("K23").Value = 0.0
$ws.Range("L23").Value = 0.0

$ws.Range("I24:L24").ClearFormats()
$ws.Range("I24").Value = 0.0
$ws.Range("J24").Value = 0.0
$ws.Range("K24").Value = 0.0
$ws.Range("L24").Value = 0.0

$ws.Range("I25:L25").ClearFormats()
$ws.Range("I25").Value = 0.0
$ws.Range("J25").Value = 0.0
$ws.Range("K25").Value = -7712.0
$ws.Range("L25").Value = 0.0

$ws.Range("K26:L26").ClearFormats()
$ws.Range("K26").Value = 0.0
$ws.Range("L26").Value = 0.0

$ws.Range("K27:L27").ClearFormats()
$ws.Range("K27").Value = 0.0
$ws.Range("L27").Value = 0.0

